$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 150 ---
$ws.Range("A150").Value = 45503.2916666667
$ws.Range("B150").Value = 0
$ws.Range("C150").Value = 2.29999995231628
$ws.Range("D150").Value = 2.29999995231628
$ws.Range("E150").Value = 2.29999995231628
$ws.Range("F150").Value = 2.29999995231628
$ws.Range("G150").Value = "'2.29999995231628"
$ws.Range("H150").Value = "LS.MI"

# --- Row 151 ---
$ws.Range("A151").Value = 45504.5691203704
$ws.Range("B151").Value = 2000
$ws.Range("C151").Value = 2.26999998092651
$ws.Range("D151").Value = 2.24000000953674
$ws.Range("E151").Value = 2.24000000953674
$ws.Range("F151").Value = 2.26999998092651
$ws.Range("G151").Value = "'2.26999998092651"
$ws.Range("H151").Value = "LS.MI"

# Reuse the existing date/time style for the new date cells (column A)
# instead of letting Excel fabricate a brand new style entry.
$ws.Range("A149").Copy()
$ws.Range("A150:A151").PasteSpecial(-4122)

# Clear the implicit "quote prefix" style added for the text-looking
# numbers in column G by re-pasting the formats from plain (unstyled)
# cells in column H.
$ws.Range("H149").Copy()
$ws.Range("G150:G151").PasteSpecial(-4122)

$wb.Save()
